$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("정리 날짜") holds dates stored as plain text like "20201026".
# Force the number format to Text before assigning so Excel keeps storing
# these as text strings (matching the original inlineStr type) instead of
# silently converting them to numeric values. Reset the style afterwards
# so no visible formatting change is introduced.
$ws.Range("B2:B4").NumberFormat = "@"

# Update file name (column A) to include the full path prefix,
# and update the date (column B) from 20201026 to 20201030
# for rows 2-4.
$ws.Range("A2").Value = "/Users/wonmyeongkwon/Desktop/Developer/LawClerk_Kwon/(19.01.29)소장.pdf"
$ws.Range("B2").Value = "20201030"

$ws.Range("A3").Value = "/Users/wonmyeongkwon/Desktop/Developer/LawClerk_Kwon/(19.03.12)답변서.pdf"
$ws.Range("B3").Value = "20201030"

$ws.Range("A4").Value = "/Users/wonmyeongkwon/Desktop/Developer/LawClerk_Kwon/(19.04.25)답변서.pdf"
$ws.Range("B4").Value = "20201030"

# Restore the default (unformatted) style on column B so the cells keep
# the same visual style as before, while the underlying value remains text.
$ws.Range("B2:B4").Style = "Normal"
